$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab) from SCD0211 to SCD0012
$ws.Name = "SCD0012"

# Update TC_ID cell from DGS-226 to SCD0012-001
$ws.Range("B2").Value = "SCD0012-001"

# Widen column B to fit the new, longer TC_ID text
$ws.Columns("B").ColumnWidth = 11.6

# Update the selected cell / view to B3 (no frozen/top-left offset anymore)
$ws.Range("B3").Select()
